$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for columns C and D (td_sim_1 and record_atd)
$updates = @{
    3  = 171
    5  = 924
    7  = 23
    9  = 166
    11 = 193
    13 = 6
    15 = 20
    17 = 14
    19 = 135
    21 = 7
    23 = 476
    24 = 324
    26 = 47
    28 = 334
    30 = 33
    32 = 51
    34 = 10
    36 = 41
    38 = 134
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("C$row").Value = $value
    $ws.Range("D$row").Value = $value
}

# Update the average formula value in C39 (td_sim_1 column average)
$ws.Range("C39").Value = 163.6315789473684
